# Update the cryptocurrency price/volume table with the latest scraped
# values. Price strings (column D) that look like plain numbers (e.g.
# "215.03") are written with a leading apostrophe so Excel stores them
# as text (matching the original inlineStr cell type) instead of
# auto-converting them to numeric values; the Style is then reset back
# to "Normal" so no residual text-formatting is left on the cell.
# Percent strings (column E) already contain surrounding whitespace so
# they are always kept as plain text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.748.94'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.633.67'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = "'215.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = "'19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '1.857.14'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '1.615.38'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = "'0.557"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = "'62.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '25.746.10'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").Value = "'193.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = "'6.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.48%  '
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  +3.55%  '
$ws.Range("D26").Value = "'142.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").Value = "'15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").Value = "'0.0492"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("E32").Value = '  +1.07%  '
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = "'0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '1.127.32'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("D39").Value = "'0.547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.60%  '
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("E42").Value = '  +2.08%  '
$ws.Range("D43").Value = "'99.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").Value = '1.767.86'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").Value = "'55.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +3.60%  '
$ws.Range("D51").Value = "'7.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.53%  '
